# Edit: insert a new row into the "Line Items" sheet and move the
# active sheet/selection state to reflect it.

$wb = $excel.ActiveWorkbook

# --- Line Items sheet: insert a new row 4 with values 1,3,5 ---
$ws = $wb.Worksheets.Item("Line Items")
[void]$ws.Activate()

[void]$ws.Rows("4:4").Insert()
$ws.Rows("4:4").RowHeight = 15.75

$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 5

[void]$ws.Range("F9").Select()

# --- Invoices sheet: move selection to A2 ---
$wsInvoices = $wb.Worksheets.Item("Invoices")
[void]$wsInvoices.Activate()
[void]$wsInvoices.Range("A2").Select()

# --- Re-activate Line Items so it is the active tab ---
[void]$ws.Activate()
